$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$block_B = New-Object 'object[,]' 24,3
$block_F = New-Object 'object[,]' 24,4
$block_K = New-Object 'object[,]' 24,3

$block_B[0,0] = 0.4775687422101953
$block_B[0,1] = 0.1281789982642039
$block_B[0,2] = 0.06245909031792962
$block_F[0,0] = 1.334852579320426
$block_F[0,1] = 1.219765546245497
$block_F[0,2] = 1.172461250130837
$block_F[0,3] = 1.234641739171472
$block_K[0,0] = 0.2901266050507445
$block_K[0,1] = 0.2986121799919061
$block_K[0,2] = 0.176813460862391
$block_B[1,0] = 0.4455515744904233
$block_B[1,1] = 0.1254283203567823
$block_B[1,2] = 0.06159383311435462
$block_F[1,0] = 1.329561281832071
$block_F[1,1] = 1.21606804412572
$block_F[1,2] = 1.175306587185418
$block_F[1,3] = 1.237561429762366
$block_K[1,0] = 0.2579919681258502
$block_K[1,1] = 0.2951446198367407
$block_K[1,2] = 0.1701167158920462
$block_B[2,0] = 0.4260965005388471
$block_B[2,1] = 0.1237085438182461
$block_B[2,2] = 0.06105191816130784
$block_F[2,0] = 1.326991855540186
$block_F[2,1] = 1.214440438525884
$block_F[2,2] = 1.177514056302286
$block_F[2,3] = 1.239879424251129
$block_K[2,0] = 0.2383123340322726
$block_K[2,1] = 0.2931621517200327
$block_K[2,2] = 0.166094313475007
$block_B[3,0] = 0.418219991795894
$block_B[3,1] = 0.1229999544659819
$block_B[3,2] = 0.06082841788668603
$block_F[3,0] = 1.326115638340987
$block_F[3,1] = 1.213938685205918
$block_F[3,2] = 1.178529455405297
$block_F[3,3] = 1.240956122254858
$block_K[3,0] = 0.2303058547777397
$block_K[3,1] = 0.2923912029010793
$block_K[3,2] = 0.1644777226173417
$block_B[4,0] = 0.4169152312571498
$block_B[4,1] = 0.122881824406555
$block_B[4,2] = 0.06079114516800743
$block_F[4,0] = 1.325980461519848
$block_F[4,1] = 1.213865121387272
$block_F[4,2] = 1.178705059761228
$block_F[4,3] = 1.241142885667159
$block_K[4,0] = 0.2289771859990424
$block_K[4,1] = 0.2922654190439999
$block_K[4,2] = 0.1642106546388646
$block_B[5,0] = 0.4259900657663422
$block_B[5,1] = 0.1236990189920704
$block_B[5,2] = 0.06104891473635021
$block_F[5,0] = 1.326979346839074
$block_F[5,1] = 1.214433017892034
$block_F[5,2] = 1.177527281230311
$block_F[5,3] = 1.239893410116807
$block_K[5,0] = 0.2382043023293647
$block_K[5,1] = 0.2931516048649101
$block_K[5,2] = 0.166072420065273
$block_B[6,0] = 0.4664871922792315
$block_B[6,1] = 0.1272369610685118
$block_B[6,2] = 0.06216296495965423
$block_F[6,0] = 1.332887109460984
$block_F[6,1] = 1.218357197777181
$block_F[6,2] = 1.173346796948863
$block_F[6,3] = 1.23553942650247
$block_K[6,0] = 0.2790361124053788
$block_K[6,1] = 0.2973861621529466
$block_K[6,2] = 0.1744859105113932
$block_B[7,0] = 0.5475044958765807
$block_B[7,1] = 0.1339307594602772
$block_B[7,2] = 0.06426277243043188
$block_F[7,0] = 1.349866067617612
$block_F[7,1] = 1.231158161411827
$block_F[7,2] = 1.168800217903907
$block_F[7,3] = 1.231169854930982
$block_K[7,0] = 0.3595055935455775
$block_K[7,1] = 0.3068520300560778
$block_K[7,2] = 0.1916916625222598
$block_B[8,0] = 0.6079942087771997
$block_B[8,1] = 0.1387012080048891
$block_B[8,2] = 0.06575338408697462
$block_F[8,0] = 1.365635835616814
$block_F[8,1] = 1.243687131038456
$block_F[8,2] = 1.167684422935992
$block_F[8,3] = 1.230503028891121
$block_K[8,0] = 0.4188654404324836
$block_K[8,1] = 0.3145140653432605
$block_K[8,2] = 0.2047617031171285
$block_B[9,0] = 0.6357205505734669
$block_B[9,1] = 0.1408396901096864
$block_B[9,2] = 0.06642011157572369
$block_F[9,0] = 1.373527195104259
$block_F[9,1] = 1.250067912893613
$block_F[9,2] = 1.167659667602564
$block_F[9,3] = 1.230752547519081
$block_K[9,0] = 0.4459212121983285
$block_K[9,1] = 0.3181532341211408
$block_K[9,2] = 0.2108004325172885
$block_B[10,0] = 0.6462495895284235
$block_B[10,1] = 0.1416449465348251
$block_B[10,2] = 0.06667094165928944
$block_F[10,0] = 1.376618713252242
$block_F[10,1] = 1.252582277087427
$block_F[10,2] = 1.167719688781077
$block_F[10,3] = 1.230926559498528
$block_K[10,0] = 0.4561739292064146
$block_K[10,1] = 0.3195533567199647
$block_K[10,2] = 0.213100469734627
$block_B[11,0] = 0.6439806594189008
$block_B[11,1] = 0.1414717219856243
$block_B[11,2] = 0.06661699426526013
$block_F[11,0] = 1.375948307628022
$block_F[11,1] = 1.252036398629585
$block_F[11,2] = 1.167703676303361
$block_F[11,3] = 1.230885545641868
$block_K[11,0] = 0.4539655038516912
$block_K[11,1] = 0.3192508354444215
$block_K[11,2] = 0.212604525356312
$block_B[12,0] = 0.6365861885965103
$block_B[12,1] = 0.1409060299993712
$block_B[12,2] = 0.06644078055522584
$block_F[12,0] = 1.373779467023965
$block_F[12,1] = 1.250272804315173
$block_F[12,2] = 1.167663214987073
$block_F[12,3] = 1.230765269659358
$block_K[12,0] = 0.4467645650923089
$block_K[12,1] = 0.3182679814499352
$block_K[12,2] = 0.2109893918771775
$block_B[13,0] = 0.6320607139223569
$block_B[13,1] = 0.1405589362284729
$block_B[13,2] = 0.0663326299782554
$block_F[13,0] = 1.37246443449277
$block_F[13,1] = 1.249205331364877
$block_F[13,2] = 1.167647467594094
$block_F[13,3] = 1.230701954318462
$block_K[13,0] = 0.4423547220622766
$block_K[13,1] = 0.317668824885331
$block_K[13,2] = 0.2100018056630901
$block_B[14,0] = 0.6061864024103158
$block_B[14,1] = 0.1385608174633575
$block_B[14,2] = 0.06570958229565349
$block_F[14,0] = 1.365134560334795
$block_F[14,1] = 1.243283852534944
$block_F[14,2] = 1.1676957535024
$block_F[14,3] = 1.230497849287993
$block_K[14,0] = 0.4170983171837577
$block_K[14,1] = 0.3142793251365532
$block_K[14,2] = 0.2043689240006969
$block_B[15,0] = 0.5903666512420784
$block_B[15,1] = 0.1373269425141501
$block_B[15,2] = 0.06532444496960466
$block_F[15,0] = 1.360821746354716
$block_F[15,1] = 1.23982580601681
$block_F[15,2] = 1.167849014045871
$block_F[15,3] = 1.230514255509725
$block_K[15,0] = 0.4016176152115918
$block_K[15,1] = 0.3122393004649098
$block_K[15,2] = 0.2009371176171086
$block_B[16,0] = 0.5812872709008445
$block_B[16,1] = 0.1366142737710447
$block_B[16,2] = 0.06510185566857274
$block_F[16,0] = 1.358408672615155
$block_F[16,1] = 1.237900947591882
$block_F[16,2] = 1.167982613277388
$block_F[16,3] = 1.230575728403039
$block_K[16,0] = 0.3927185006942011
$block_K[16,1] = 0.3110803993708089
$block_K[16,2] = 0.1989720017455454
$block_B[17,0] = 0.5782165487687791
$block_B[17,1] = 0.1363724649694262
$block_B[17,2] = 0.06502630760032702
$block_F[17,0] = 1.357603246630347
$block_F[17,1] = 1.237260231986383
$block_F[17,2] = 1.168035654441979
$block_F[17,3] = 1.230605478656464
$block_K[17,0] = 0.3897062769445938
$block_K[17,1] = 0.310690501619618
$block_K[17,2] = 0.1983081558988715
$block_B[18,0] = 0.5920486526269144
$block_B[18,1] = 0.1374585986196308
$block_B[18,2] = 0.0653655541656093
$block_F[18,0] = 1.361273862046801
$block_F[18,1] = 1.240187283990608
$block_F[18,2] = 1.16782799576626
$block_F[18,3] = 1.230507123283594
$block_K[18,0] = 0.4032650495070698
$block_K[18,1] = 0.3124549676364978
$block_K[18,2] = 0.2013015323978777
$block_B[19,0] = 0.6387573218651141
$block_B[19,1] = 0.1410723106199043
$block_B[19,2] = 0.06649258353450449
$block_F[19,0] = 1.374413706094884
$block_F[19,1] = 1.250788151246965
$block_F[19,2] = 1.167673216353919
$block_F[19,3] = 1.230798439163038
$block_K[19,0] = 0.4488794600071628
$block_K[19,1] = 0.3185560714994011
$block_K[19,2] = 0.2114634354221678
$block_B[20,0] = 0.6694569140685189
$block_B[20,1] = 0.143407618781751
$block_B[20,2] = 0.06721956789591133
$block_F[20,0] = 1.383603084263342
$block_F[20,1] = 1.258288296719343
$block_F[20,2] = 1.167976533459424
$block_F[20,3] = 1.231452372771855
$block_K[20,0] = 0.4787334332327191
$block_K[20,1] = 0.3226719874497803
$block_K[20,2] = 0.2181823175732944
$block_B[21,0] = 0.653056292236414
$block_B[21,1] = 0.1421636399580848
$block_B[21,2] = 0.0668324443783419
$block_F[21,0] = 1.37864346889873
$block_F[21,1] = 1.254232958818065
$block_F[21,2] = 1.167777649186846
$block_F[21,3] = 1.231060934094067
$block_K[21,0] = 0.4627960296848528
$block_K[21,1] = 0.3204635049012126
$block_K[21,2] = 0.2145892617998015
$block_B[22,0] = 0.5912881711623754
$block_B[22,1] = 0.1373990871669548
$block_B[22,2] = 0.06534697233723108
$block_F[22,0] = 1.36106925364588
$block_F[22,1] = 1.240023662883004
$block_F[22,2] = 1.16783735643078
$block_F[22,3] = 1.230510185656506
$block_K[22,0] = 0.4025202414921409
$block_K[22,1] = 0.3123574210954558
$block_K[22,2] = 0.2011367559388049
$block_B[23,0] = 0.5254167146489124
$block_B[23,1] = 0.1321459053769232
$block_B[23,2] = 0.06370384328375067
$block_F[23,0] = 1.34469462530879
$block_F[23,1] = 1.22714736890353
$block_F[23,2] = 1.169639431245855
$block_F[23,3] = 1.231905435560869
$block_K[23,0] = 0.3376942299806842
$block_K[23,1] = 0.3041669614807176
$block_K[23,2] = 0.1869615762871106

$ws.Range("B2:D25").Value = $block_B
$ws.Range("F2:I25").Value = $block_F
$ws.Range("K2:M25").Value = $block_K
